$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the columns we touch are formatted as text so that numeric-looking
# strings (e.g. "566.89", "0.0584") are preserved exactly as text, matching
# the original inline-string cell content instead of being coerced to numbers.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "63.850.71"
$ws.Range("E2").Value = "  +0.09%  "
$ws.Range("D3").Value = "2.736.05"
$ws.Range("E3").Value = "  -0.52%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "566.89"
$ws.Range("E5").Value = "  -0.81%  "
$ws.Range("D6").Value = "161.31"
$ws.Range("E6").Value = "  +2.71%  "
$ws.Range("D7").Value = "0.998"
$ws.Range("E7").Value = "  -0.07%  "
$ws.Range("D8").Value = "0.598"
$ws.Range("E8").Value = "  -0.22%  "
$ws.Range("D9").Value = "0.110"
$ws.Range("E9").Value = "  +0.79%  "
$ws.Range("E10").Value = "  +4.35%  "
$ws.Range("D11").Value = "5.67"
$ws.Range("E11").Value = "  +0.54%  "
$ws.Range("D12").Value = "0.381"
$ws.Range("E12").Value = "  +0.05%  "
$ws.Range("D13").Value = "3.218.13"
$ws.Range("E13").Value = "  -0.57%  "
$ws.Range("D14").Value = "26.99"
$ws.Range("E14").Value = "  +2.31%  "
$ws.Range("D15").Value = "63.667.98"
$ws.Range("E15").Value = "  +0.34%  "
$ws.Range("E16").Value = "  +0.65%  "
$ws.Range("D17").Value = "2.735.66"
$ws.Range("E17").Value = "  -0.68%  "
$ws.Range("D18").Value = "12.64"
$ws.Range("E18").Value = "  +4.60%  "
$ws.Range("E19").Value = "  -0.74%  "
$ws.Range("D20").Value = "355.61"
$ws.Range("E20").Value = "  +0.54%  "
$ws.Range("D21").Value = "6.61"
$ws.Range("E21").Value = "  -1.48%  "
$ws.Range("D22").Value = "0.998"
$ws.Range("E22").Value = "  +0.06%  "
$ws.Range("D23").Value = "0.523"
$ws.Range("E23").Value = "  -2.22%  "
$ws.Range("D24").Value = "64.68"
$ws.Range("E24").Value = "  -0.68%  "
$ws.Range("D25").Value = "0.170"
$ws.Range("E25").Value = "  +0.55%  "
$ws.Range("D26").Value = "0.999"
$ws.Range("E26").Value = "  -0.09%  "
$ws.Range("D27").Value = "8.42"
$ws.Range("E27").Value = "  +0.51%  "
$ws.Range("D28").Value = "0.0₃0914"
$ws.Range("E28").Value = "  +1.76%  "
$ws.Range("D29").Value = "1.99"
$ws.Range("E29").Value = "  +3.46%  "
$ws.Range("D30").Value = "7.19"
$ws.Range("E30").Value = "  +3.34%  "
$ws.Range("D31").Value = "1.34"
$ws.Range("E31").Value = "  +11.88%  "
$ws.Range("D32").Value = "166.55"
$ws.Range("E32").Value = "  -1.43%  "
$ws.Range("D33").Value = "4.95"
$ws.Range("E33").Value = "  +2.34%  "
$ws.Range("D34").Value = "20.13"
$ws.Range("E34").Value = "  +0.14%  "
$ws.Range("D35").Value = "1.49"
$ws.Range("E35").Value = "  +4.12%  "
$ws.Range("E36").Value = "  -0.03%  "
$ws.Range("D37").Value = "1.82"
$ws.Range("E37").Value = "  +2.08%  "
$ws.Range("D38").Value = "0.981"
$ws.Range("E38").Value = "  +0.63%  "
$ws.Range("D39").Value = "348.94"
$ws.Range("E39").Value = "  +7.03%  "
$ws.Range("D40").Value = "6.36"
$ws.Range("E40").Value = "  +3.50%  "
$ws.Range("D41").Value = "4.10"
$ws.Range("E41").Value = "  -0.61%  "
$ws.Range("D42").Value = "38.70"
$ws.Range("E42").Value = "  -0.47%  "
$ws.Range("D43").Value = "21.98"
$ws.Range("E43").Value = "  +3.63%  "
$ws.Range("D44").Value = "21.19"
$ws.Range("E44").Value = "  -0.61%  "
$ws.Range("D45").Value = "0.0587"
$ws.Range("E45").Value = "  +0.98%  "
$ws.Range("D46").Value = "0.627"
$ws.Range("E46").Value = "  +0.77%  "
$ws.Range("D47").Value = "0.0252"
$ws.Range("E47").Value = "  -0.08%  "
$ws.Range("E48").Value = "  -0.10%  "
$ws.Range("D51").Value = "11.08"
$ws.Range("E51").Value = "  +0.30%  "

# Rows 49 and 50 swap order (Aave <-> FirstDigitalUSD) plus value updates
$ws.Range("B49").Value = "Aave"
$ws.Range("C49").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D49").Value = "132.42"
$ws.Range("E49").Value = "  -1.83%  "
$ws.Range("B50").Value = "FirstDigitalUSD"
$ws.Range("C50").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D50").Value = "0.998"
$ws.Range("E50").Value = "  -0.12%  "
